$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (old D25/J25 values) is inserted as a new row 26,
# shifting the two most recent weekly entries down by one row:
#   old row24 -> new row25
#   old row25 -> new row26
#   row24 gets the brand-new weekly values

# Copy formatting of the last data row (25) down into the new row 26
$ws.Range("A25:R25").Copy() | Out-Null
$ws.Range("A26:R26").PasteSpecial(-4122) | Out-Null

# New row 26 = old row 25 values
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = "Terminal La Palmera de La Serena"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 45007
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112039
$ws.Range("G26").Value = "Ciboulette"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 1160
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2500
$ws.Range("M26").Value = 2250
$ws.Range("N26").Value = "`$/docena de atados"
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 750
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = "Hortaliza"

# Row 25 becomes what used to be row 24 (only date & volume differ)
$ws.Range("D25").Value = 44970
$ws.Range("J25").Value = 800

# Row 24 receives the brand-new weekly values
$ws.Range("D24").Value = 45041
$ws.Range("J24").Value = 1160
